# household_member.xlsx revision:
#  - survey sheet: insert 2 new field rows (household_id, member_name) at top of data
#  - settings sheet: add a new "table_id"/"household_member" setting row
#  - add new "model" sheet describing a schema.joins example

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "survey" sheet (sheet1) -- insert two rows after the header row
# ---------------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")

# Insert two blank rows at row 2 (pushes the existing data rows 2-7 to 4-9)
$survey.Rows.Item(2).Resize(2).Insert()

# Row 2: household_id (readonly text field with hint + comment + readonly flag)
$survey.Range("A2").Value = "text"
$survey.Range("B2").Value = "household_id"
$survey.Range("C2").Value = "Unique barcode ID or locator designation for this household (foreign key into household table)."
$survey.Range("D2").Value = $null
$survey.Range("E2").Value = $true
$survey.Range("F2").Value = "readonly is not implemented, but this would be a note-like field that was set to be the parent table key via the extra URL arguments.  It could be hidden and only displayed as a reference in a prompt label."

# Row 3: member_name
$survey.Range("A3").Value = "text"
$survey.Range("B3").Value = "member_name"
$survey.Range("C3").Value = "Enter the name of the household member:"
$survey.Range("D3").Value = $null

# Header cells for the new "readonly" / "comments" columns (row 1)
$survey.Range("E1").Value = "readonly"
$survey.Range("F1").Value = "comments"

# Formatting to match the new field rows: font applied (no family-2 arial for
# col A/B, with family-2 arial + wrap for col C), row 2 tall enough to show
# the long comment, column F widened for the comment text.
$survey.Range("A2:B3").Font.Name = "Arial"
$survey.Range("A2:B3").Font.Size = 10
$survey.Range("A2:B3").WrapText = $false

$survey.Range("C2:C3").Font.Name = "Arial"
$survey.Range("C2:C3").Font.Size = 10
$survey.Range("C2:C3").WrapText = $true

$survey.Rows.Item(2).RowHeight = 103.8
$survey.Rows.Item(3).RowHeight = 13.2
$survey.Columns.Item(6).ColumnWidth = 25.6640625

$survey.Range("C9").Select()

# ---------------------------------------------------------------------------
# 2. "settings" sheet (sheet3) -- append a table_id / household_member row
# ---------------------------------------------------------------------------
$settings = $wb.Worksheets.Item("settings")

$settings.Range("A5").Value = "table_id"
$settings.Range("B5").Value = "household_member"
$settings.Range("A5:B5").Font.Name = "Arial"
$settings.Range("A5:B5").Font.Size = 10
$settings.Range("A5:B5").Font.Bold = $false
$settings.Range("A5:B5").Font.Color = 0
$settings.Range("A5:B5").WrapText = $true

$settings.Range("B6").Select()

# ---------------------------------------------------------------------------
# 3. New "model" sheet describing an example schema.joins column_definition
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$model = $wb.Worksheets.Add($null, $lastSheet)
$model.Name = "model"

$model.Range("A1").Value = "name"
$model.Range("B1").Value = "schema.type"
$model.Range("C1").Value = "schema.joins[0].table_id"
$model.Range("D1").Value = "schema.joins[0].element_name"
$model.Range("E1").Value = "comment"

$model.Range("A2").Value = "household_id"
$model.Range("B2").Value = "string"
$model.Range("C2").Value = "household"
$model.Range("D2").Value = "household_id"

# Rich-text comment cell with bold "household_id" / "household_member" runs.
# Only the two bold runs are touched -- the plain-text runs around them are
# left alone so they keep the workbook's implicit default formatting (no
# explicit run properties), matching how the surrounding text was authored.
$commentCell = $model.Range("E2")
$commentText = "This would insert a 'joins' entry into the column_definitions table for the household_id column of the household_member table_id of the form: " + `
    [char]10 + '"[ { table_id: household, element_name: household_id } ]"' + [char]10 + `
    "The way to define joins is undoubtedly broken in the current XLSXConverter, as there is no way to parse lists of values (as far as I know). I will research how to fix this."
$commentCell.Value = $commentText

$run1Len = ("This would insert a 'joins' entry into the column_definitions table for the ").Length

$run2Start = $run1Len + 1
$run2Len = ("household_id").Length
$commentCell.Characters($run2Start, $run2Len).Font.Bold = $true
$commentCell.Characters($run2Start, $run2Len).Font.Name = "Arial"
$commentCell.Characters($run2Start, $run2Len).Font.Size = 10

$run3Len = (" column of the ").Length
$run4Start = $run2Start + $run2Len + $run3Len
$run4Len = ("household_member").Length
$commentCell.Characters($run4Start, $run4Len).Font.Bold = $true
$commentCell.Characters($run4Start, $run4Len).Font.Name = "Arial"
$commentCell.Characters($run4Start, $run4Len).Font.Size = 10

$commentCell.WrapText = $true

$model.Rows.Item(2).RowHeight = 225
$model.Columns.Item("A:D").ColumnWidth = 34.33203125
$model.Columns.Item("E").ColumnWidth = 46.109375

$model.PageSetup.Orientation = 1

$model.Range("C12").Select()
